$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: assign directly. ---
$ws.Range("D2").Value = "42.851.15"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.533.64"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "2.921.09"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.577.67"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").Value = "42.832.23"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E33").Value = "  +8.63%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "1.988.81"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D48").Value = "2.772.96"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -2.57%  "

# --- Numeric-looking strings must stay text (source values like
# "42.851.15" or "311.44" are plain text in the sheet, not numbers).
# Writing such a string straight into .Value lets Excel auto-convert it
# to a Double (dropping e.g. a trailing zero: "69.70" -> 69.7). Instead,
# put a text-literal formula in the cell, then Copy + PasteSpecial
# (values only) back onto itself so the stored result is a plain text
# cell (no formula, no stray number format) - done one cell at a time so
# the paste lands back on the exact source cell instead of being
# compacted across a multi-area selection.
$ws.Range("D5").Formula = '="311.44"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="101.05"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D7").Formula = '="0.566"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="0.524"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = '="35.84"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.0807"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="7.37"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D16").Formula = '="15.42"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D17").Formula = '="0.817"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="6.72"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="12.36"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D22").Formula = '="69.70"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D23").Formula = '="244.17"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="2.89"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="2.04"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="25.56"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D29").Formula = '="10.20"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="38.89"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = '="162.07"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="5.82"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D34").Formula = '="2.67"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D35").Formula = '="0.0791"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="18.43"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="0.111"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D40").Formula = '="0.118"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="4.21"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="22.35"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="3.32"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="0.0301"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D47").Formula = '="9.34"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="0.192"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("D50").Formula = '="79.64"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="72.43"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
